$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "25.911.03"
$ws.Cells.Item(2,5).Value = "  -0.50%  "
$ws.Cells.Item(3,4).Value = "1.742.49"
$ws.Cells.Item(3,5).Value = "  -1.13%  "
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = "1.000"
$ws.Cells.Item(4,5).Value = "  -0.19%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "231.11"
$ws.Cells.Item(5,5).Value = "  -2.82%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "1.001"
$ws.Cells.Item(6,5).Value = "  -0.11%  "
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.5249"
$ws.Cells.Item(7,5).Value = "  +0.39%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "0.2759"
$ws.Cells.Item(8,5).Value = "  +0.07%  "
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = "39.54"
$ws.Cells.Item(9,5).Value = "  -2.26%  "
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.06142"
$ws.Cells.Item(10,5).Value = "  -0.97%  "
$ws.Cells.Item(11,4).Value = "1.736.82"
$ws.Cells.Item(11,5).Value = "  -1.59%  "
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.07118"
$ws.Cells.Item(12,5).Value = "  +1.34%  "
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "15.23"
$ws.Cells.Item(13,5).Value = "  -3.04%  "
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "0.6454"
$ws.Cells.Item(14,5).Value = "  +0.65%  "
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "4.523"
$ws.Cells.Item(15,5).Value = "  -0.53%  "
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "77.26"
$ws.Cells.Item(16,5).Value = "  -1.18%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "1.000"
$ws.Cells.Item(17,5).Value = "  -0.14%  "
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "1.0000"
$ws.Cells.Item(18,5).Value = "  -0.14%  "
$ws.Cells.Item(19,4).Value = "25.897.78"
$ws.Cells.Item(19,5).Value = "  -0.64%  "
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "11.55"
$ws.Cells.Item(20,5).Value = "  -0.95%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.000006674"
$ws.Cells.Item(21,5).Value = "  -0.92%  "
$ws.Cells.Item(22,4).Value = "1.959.63"
$ws.Cells.Item(22,5).Value = "  -2.17%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "4.275"
$ws.Cells.Item(23,5).Value = "  +4.85%  "
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "8.783"
$ws.Cells.Item(24,5).Value = "  +3.89%  "
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "5.182"
$ws.Cells.Item(25,5).Value = "  -0.32%  "
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "139.91"
$ws.Cells.Item(26,5).Value = "  +0.81%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "1.522"
$ws.Cells.Item(27,5).Value = "  +0.63%  "
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "15.20"
$ws.Cells.Item(28,5).Value = "  +0.06%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "1.803"
$ws.Cells.Item(29,5).Value = "  -2.04%  "
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "102.53"
$ws.Cells.Item(30,5).Value = "  -0.88%  "
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0.08327"
$ws.Cells.Item(31,5).Value = "  -0.75%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "3.733"
$ws.Cells.Item(32,5).Value = "  +0.90%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "3.577"
$ws.Cells.Item(33,5).Value = "  +3.67%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "0.04525"
$ws.Cells.Item(34,5).Value = "  +1.28%  "
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "2.613"
$ws.Cells.Item(35,5).Value = "  -0.43%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.9791"
$ws.Cells.Item(36,5).Value = "  -2.35%  "
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "0.6229"
$ws.Cells.Item(37,5).Value = "  +2.74%  "
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "2.701"
$ws.Cells.Item(38,5).Value = "  -1.39%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.01590"
$ws.Cells.Item(39,5).Value = "  -0.11%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "1.927"
$ws.Cells.Item(40,5).Value = "  -3.07%  "
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "1.000"
$ws.Cells.Item(41,5).Value = "  -0.16%  "
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "100.22"
$ws.Cells.Item(42,5).Value = "  -2.48%  "
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.3875"
$ws.Cells.Item(43,5).Value = "  -0.11%  "
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "0.7337"
$ws.Cells.Item(44,5).Value = "  -1.50%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "5.023"
$ws.Cells.Item(45,5).Value = "  +1.70%  "
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.05339"
$ws.Cells.Item(46,5).Value = "  -3.22%  "
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.1125"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "6.260"
$ws.Cells.Item(48,5).Value = "  -1.27%  "
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "53.68"
$ws.Cells.Item(49,5).Value = "  +1.90%  "
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "30.20"
$ws.Cells.Item(50,5).Value = "  -0.06%  "
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "7.648"
$ws.Cells.Item(51,5).Value = "  +2.89%  "
